# Apply the new ("Office Theme") colour palette to the presentation's
# design / theme so that every slide (via the slide master) picks up the
# Office Theme colours instead of the previous "Integral" theme colours.
#
# PowerPoint stores a theme's 12-colour palette (dk1, lt1, dk2, lt2,
# accent1-6, hlink, folHlink) on the Theme object that hangs off the
# slide master. We can't "import" a whole external .thmx/theme file via
# COM automation, but we *can* repaint each of the twelve colour slots
# individually through ThemeColorScheme.Colors(i).RGB - which is exactly
# how PowerPoint itself represents "switch to a different theme" at the
# object-model level when no new master/layout structure is required.

$p = $ppt.ActivePresentation
$master = $p.SlideMaster
$tcs = $master.Theme.ThemeColorScheme

# index -> (scheme slot, new "Office Theme" RGB as BGR-packed long for
# the COM RGB() encoding used by PowerPoint: val = R + G*256 + B*65536)
$tcs.Colors(1).RGB  = 0        # dk1      000000
$tcs.Colors(2).RGB  = 16777215 # lt1      FFFFFF
$tcs.Colors(3).RGB  = 6968388  # dk2      44546A
$tcs.Colors(4).RGB  = 15132391 # lt2      E7E6E6
$tcs.Colors(5).RGB  = 13998939 # accent1  5B9BD5
$tcs.Colors(6).RGB  = 3243501  # accent2  ED7D31
$tcs.Colors(7).RGB  = 10855845 # accent3  A5A5A5
$tcs.Colors(8).RGB  = 49407    # accent4  FFC000
$tcs.Colors(9).RGB  = 12874308 # accent5  4472C4
$tcs.Colors(10).RGB = 4697456  # accent6  70AD47
$tcs.Colors(11).RGB = 12673797 # hlink    0563C1
$tcs.Colors(12).RGB = 7491477  # folHlink 954F72
